# Add a new question/answer pair (row 3) to the Questions sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "What is the address of the University of California, Berkeley?"
$ws.Range("B3").Value = "SELECT pcaddr`nFROM ic2022campuses`nWHERE pcinstnm = 'University of California, Berkeley';"

# Match the wrapped-text style used by the other "Expected SQL Query" cells.
$ws.Range("B3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 85

# Move the active selection, as it was left after the edit.
$ws.Range("C3").Select()
